# Fruta / hortaliza, semanal
# Weekly refresh: insert the new week's price record at row 3 (pushing the
# older history rows down by one), matching the rest of the sheet's
# template (market / product / region columns) and only varying the
# date, volume and price columns for the new observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 3; rows 3..11 shift down to 4..12.
$ws.Rows.Item(3).Insert()

# Fill the new row 3 with this week's record.
$ws.Range("A3").Value = 10
$ws.Range("B3").Value = "Vega Modelo de Temuco"
$ws.Range("C3").Value = "La Araucanía"
$ws.Range("D3").Value = 44998
$ws.Range("E3").Value = 9
$ws.Range("F3").Value = "Fruta"
$ws.Range("G3").Value = 100101
$ws.Range("H3").Value = "Berries"
$ws.Range("I3").Value = 100101004
$ws.Range("J3").Value = "Frambuesa"
$ws.Range("K3").Value = "Sin especificar"
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 20
$ws.Range("N3").Value = 2500
$ws.Range("O3").Value = 2500
$ws.Range("P3").Value = 2500
$ws.Range("Q3").Value = "$/envase 1 kilo"
$ws.Range("R3").Value = "Región de La Araucanía"
$ws.Range("S3").Value = 2500
$ws.Range("T3").Value = 1
